$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data / mean calculation
$ws.Range("F3").Value = -9
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -8
$ws.Range("F6").Value = -4
$ws.Range("F12").Value = -3
$ws.Range("F14").Value = 5
